$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Extend the table with a new year column (2023), copying the formatting
# from the previous year's column (J) so styles/number formats match.
$ws.Range("J3:J6").Copy($ws.Range("K3:K6"))

$ws.Range("K3").Value = 2023
$ws.Range("K4").Value = 1425.3
$ws.Range("K5").Value = 859.5
$ws.Range("K6").Value = 1642.2

$wb.Save()
